$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCall")

# Insert two new columns at L:M (existing L,M,N.. shift right to N,O,P..)
$ws.Range("L1:M1").EntireColumn.Insert() | Out-Null

# Header row: new "Send Payment Notification" / "Send Call Notice" columns.
# Give them the same header style as the other header cells by copying K1's format first.
$ws.Range("K1").Copy($ws.Range("L1"))
$ws.Range("K1").Copy($ws.Range("M1"))
$ws.Range("M1").Value() = "Send Call Notice"
$ws.Range("L1").Value() = "Send Payment Notification"

# Data rows: default "Yes" for the new columns.
$ws.Range("L2").Value() = "Yes"
$ws.Range("M2").Value() = "Yes"
$ws.Range("L3").Value() = "Yes"
$ws.Range("M3").Value() = "Yes"
$ws.Range("L4").Value() = "Yes"
$ws.Range("M4").Value() = "Yes"

# Column widths: L & M match K's width. The former N column (now P) keeps
# its pre-existing width/bestFit automatically via the column insert/shift.
$ws.Columns.Item(12).ColumnWidth = 20.71
$ws.Columns.Item(13).ColumnWidth = 20.71

# Data validation: extend K6:K1048576 list validation to also cover L:M.
$ws.Range("K6:K1048576").Validation.Delete()
$ws.Range("K6:M1048576").Validation.Add(3, 1, 1, """Percentage of Commitment,Upload""")

# Selection moves to A4.
$ws.Range("A4").Select() | Out-Null
